$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("US in Tasks zerlegt")

# ---------------------------------------------------------------------
# Section "1. Userstory" (rows 6-10): mark the first two tasks done and
# fill in the third task's (previously empty) responsible/deadline/status.
# ---------------------------------------------------------------------
$ws.Range("L8").Value = "abgeschlossen"
$ws.Range("J9").Value = "Simon"
$ws.Range("L9").Value = "abgeschlossen"

$ws.Range("J10").Value = "Manuel"
$ws.Range("K10").Value = (Get-Date -Year 2019 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K10").NumberFormat = $ws.Range("K9").NumberFormat
$ws.Range("L10").Value = "abgeschlossen"

# ---------------------------------------------------------------------
# Section "2. Userstory" used to be just a plain list of other user
# stories (rows 14-18). It now gets its own "broken into tasks" block,
# just like the first userstory, and the leftover userstory list is
# pushed out to columns L:R (rows 16-19).
# ---------------------------------------------------------------------

# Free up / reset the cells we are about to repurpose - they were each
# part of a C:I merge + centered style that only made sense for the old
# "plain listing" layout.
$ws.Range("C14:I14").UnMerge()
$ws.Range("C15:I15").UnMerge()
$ws.Range("C16:I16").UnMerge()
$ws.Range("C17:I17").UnMerge()
$ws.Range("C18:I18").UnMerge()
$ws.Range("C14:I18").Clear()

# New userstory header row (mirrors row 7 for the 1st userstory).
$ws.Range("B13").Value = "2. Userstory"
$ws.Range("C13").Value = "Als Guest soll man an eine Anmeldeseite kommen um sich anzumelden"
$ws.Range("C13:I13").HorizontalAlignment = $ws.Range("C7").HorizontalAlignment
$ws.Range("C13:I13").Merge()

# Task 1
$ws.Range("C14").Value = "1. Task"
$ws.Range("D14").Value = "Datenbank erzeugen"
$ws.Range("C14:I14").HorizontalAlignment = -4131
$ws.Range("J14").Value = "Simon"
$ws.Range("K14").Value = (Get-Date -Year 2019 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K14").NumberFormat = $ws.Range("K9").NumberFormat

# Task 2
$ws.Range("C15").Value = "2. Task"
$ws.Range("D15").Value = "Datenbankmethoden schreiben "
$ws.Range("J15").Value = "Simon"
$ws.Range("K15").Value = (Get-Date -Year 2019 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K15").NumberFormat = $ws.Range("K9").NumberFormat

# Task 3
$ws.Range("C16").Value = "3. Task "
$ws.Range("D16").Value = "Die Methode zum Registrieren schreiben"
$ws.Range("J16").Value = "Simon"
$ws.Range("K16").Value = (Get-Date -Year 2019 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K16").NumberFormat = $ws.Range("K9").NumberFormat

# Task 4
$ws.Range("C17").Value = "4. Task"
$ws.Range("D17").Value = "Die Methode zum Anmelden schreiben"
$ws.Range("J17").Value = "Simon"
$ws.Range("K17").Value = (Get-Date -Year 2019 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K17").NumberFormat = $ws.Range("K9").NumberFormat

# The userstories that used to be listed under column C (rows 14-18) are
# now displayed to the right of the new task table, in columns L:R.
$ws.Range("L16").Value = "Als Admin soll man die Möglichkeit haben die registrierten Benutzer verwalten können"
$ws.Range("L17").Value = "Als Admin soll man die Möglichkeit haben das Layout der Seite zu ändern"
$ws.Range("L18").Value = "Als registrierter Benutzer soll man seine Profil bearbbeiten können "
$ws.Range("L19").Value = "Als registrierter Benutzer soll man mit anderen kommunizieren können verfassen können "

$ws.Range("L16:R16").HorizontalAlignment = $ws.Range("C7").HorizontalAlignment
$ws.Range("L17:R17").HorizontalAlignment = $ws.Range("C7").HorizontalAlignment
$ws.Range("L18:R18").HorizontalAlignment = $ws.Range("C7").HorizontalAlignment
$ws.Range("L19:R19").HorizontalAlignment = $ws.Range("C7").HorizontalAlignment

$ws.Range("L16:R16").Merge()
$ws.Range("L17:R17").Merge()
$ws.Range("L18:R18").Merge()
$ws.Range("L19:R19").Merge()

$ws.Range("J23").Select()
